# Update the grid_cell (AG) column on the "solar" sheet, rows 4-26.
# These values got reshuffled among the 23 rows in the source model export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$newGridCell = @{
    4  = "CHE_3"
    5  = "CHE_12"
    6  = "CHE_0"
    7  = "CHE_17"
    8  = "CHE_19"
    9  = "CHE_14"
    10 = "CHE_18"
    11 = "CHE_10"
    12 = "CHE_22"
    13 = "CHE_13"
    14 = "CHE_20"
    15 = "CHE_1"
    16 = "CHE_6"
    17 = "CHE_7"
    18 = "CHE_11"
    19 = "CHE_15"
    20 = "CHE_25"
    21 = "CHE_24"
    22 = "CHE_5"
    23 = "CHE_8"
    24 = "CHE_21"
    25 = "CHE_9"
    26 = "CHE_4"
}

foreach ($row in $newGridCell.Keys) {
    $ws.Cells.Item($row, 33).Value = $newGridCell[$row]
}
